# Adds in the data for 2/4/22: Wordle #230, "PLEAT"
#  - Words sheet: new row 14 (date, id, word)
#  - Results sheet: six new score rows (62-67), one per player
#  - Results sheet: extends the "formatted but empty" template rows a
#    little further (68-71 pick up the A/B/C look of the rows above them)

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsWords   = $wb.Worksheets.Item("Words")

# ---------------------------------------------------------------
# Words sheet: add row 14 => 2/4/2022, Wordle 230, "PLEAT"
# ---------------------------------------------------------------
$wsWords.Cells.Item(14,1).Formula = "=A13+1"
$wsWords.Cells.Item(14,2).Formula = "=B13+1"
$wsWords.Cells.Item(14,3).Value = "PLEAT"

# Match the look of the row above it.
$wsWords.Range("A13:C13").Copy()
$wsWords.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Recalculate so the named ranges (wordle_ids/wordles/wordle_dates), which
# are based on COUNTA, pick up the new Words row before we rely on them.
$excel.CalculateFull()

# ---------------------------------------------------------------
# Results sheet: six new rows (one per player) for Wordle 230
# ---------------------------------------------------------------
$players = @(
    @{ Row = 62; Attempts = 6 },
    @{ Row = 63; Attempts = 5 },
    @{ Row = 64; Attempts = 3 },
    @{ Row = 65; Attempts = 5 },
    @{ Row = 66; Attempts = 5 },
    @{ Row = 67; Attempts = 4 }
)

$wsResults.Range("A61:C61").Copy()
$wsResults.Range("A62:C67").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($p in $players) {
    $r = $p.Row
    $prev = $r - 6

    $wsResults.Cells.Item($r,1).Formula = "=A$prev+1"
    $wsResults.Cells.Item($r,2).Formula = "=B$prev+1"
    $wsResults.Cells.Item($r,3).Formula = "=C$prev"
    $wsResults.Cells.Item($r,5).Value = $p.Attempts
}

$excel.CalculateFull()

foreach ($p in $players) {
    $r = $p.Row
    $wsResults.Cells.Item($r,4).FormulaArray = '=IF(XLOOKUP($B' + $r + ',wordle_ids,wordles)="","",XLOOKUP($B' + $r + ',wordle_ids,wordles))'
}

$excel.CalculateFull()

# ---------------------------------------------------------------
# Extend the formatted-but-empty rows a bit further (68-71), matching
# the look of the rows immediately above them.
# ---------------------------------------------------------------
$wsResults.Range("A61:C61").Copy()
$wsResults.Range("A68:C71").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# View state: active cell / scroll position on each sheet
# ---------------------------------------------------------------
$wsWords.Activate()
$wsWords.Range("C15").Select()

$wsResults.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$wsResults.Range("E66").Select()
